$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.43094266666667
$ws.Range("H2").Value = 52.292828
$ws.Range("I2").Value = 0.3741631566231533
$ws.Range("J2").Value = 0.3741631566231532
$ws.Range("M2").Value = 121.7032956666666
$ws.Range("N2").Value = 365.109887
$ws.Range("O2").Value = 0.8275728186333362
$ws.Range("P2").Value = 0.8275728186333362
$ws.Range("Q2").Value = 2121.403169110048
$ws.Range("R2").Value = 19092.62852199043
$ws.Range("S2").Value = 0.3096472581553694
$ws.Range("T2").Value = 0.3096472581553694

$ws.Range("G3").Value = 17.43094266666667
$ws.Range("H3").Value = 52.292828
$ws.Range("I3").Value = 0.3741631566231533
$ws.Range("J3").Value = 0.3741631566231532
$ws.Range("O3").Value = 0.1584547531793621
$ws.Range("P3").Value = 0.1584547531793621
$ws.Range("Q3").Value = 406.1834898231258
$ws.Range("R3").Value = 3655.651408408132
$ws.Range("S3").Value = 0.05928793063153277
$ws.Range("T3").Value = 0.05928793063153275

$ws.Range("G4").Value = 17.43094266666667
$ws.Range("H4").Value = 52.292828
$ws.Range("I4").Value = 0.3741631566231533
$ws.Range("J4").Value = 0.3741631566231532
$ws.Range("O4").Value = 0.01397242818730167
$ws.Range("P4").Value = 0.01397242818730167
$ws.Range("Q4").Value = 35.81697316455378
$ws.Range("R4").Value = 322.352758480984
$ws.Range("S4").Value = 0.005227967836251119
$ws.Range("T4").Value = 0.005227967836251117

$ws.Range("I5").Value = 0.2007058543586177
$ws.Range("J5").Value = 0.2007058543586177
$ws.Range("M5").Value = 121.7032956666666
$ws.Range("N5").Value = 365.109887
$ws.Range("O5").Value = 0.8275728186333362
$ws.Range("P5").Value = 0.8275728186333362
$ws.Range("Q5").Value = 1137.947518237727
$ws.Range("R5").Value = 10241.52766413954
$ws.Range("S5").Value = 0.1660987096077731
$ws.Range("T5").Value = 0.1660987096077731

$ws.Range("I6").Value = 0.2007058543586177
$ws.Range("J6").Value = 0.2007058543586177
$ws.Range("O6").Value = 0.1584547531793621
$ws.Range("P6").Value = 0.1584547531793621
$ws.Range("S6").Value = 0.03180279661404776
$ws.Range("T6").Value = 0.03180279661404775

$ws.Range("I7").Value = 0.2007058543586177
$ws.Range("J7").Value = 0.2007058543586177
$ws.Range("O7").Value = 0.01397242818730167
$ws.Range("P7").Value = 0.01397242818730167
$ws.Range("S7").Value = 0.002804348136796814
$ws.Range("T7").Value = 0.002804348136796813

$ws.Range("I8").Value = 0.4251309890182291
$ws.Range("J8").Value = 0.425130989018229
$ws.Range("M8").Value = 121.7032956666666
$ws.Range("N8").Value = 365.109887
$ws.Range("O8").Value = 0.8275728186333362
$ws.Range("P8").Value = 0.8275728186333362
$ws.Range("Q8").Value = 2410.376894212764
$ws.Range("R8").Value = 21693.39204791488
$ws.Range("S8").Value = 0.3518268508701938
$ws.Range("T8").Value = 0.3518268508701937

$ws.Range("I9").Value = 0.4251309890182291
$ws.Range("J9").Value = 0.425130989018229
$ws.Range("O9").Value = 0.1584547531793621
$ws.Range("P9").Value = 0.1584547531793621
$ws.Range("S9").Value = 0.0673640259337816
$ws.Range("T9").Value = 0.06736402593378157

$ws.Range("I10").Value = 0.4251309890182291
$ws.Range("J10").Value = 0.425130989018229
$ws.Range("O10").Value = 0.01397242818730167
$ws.Range("P10").Value = 0.01397242818730167
$ws.Range("S10").Value = 0.005940112214253743
$ws.Range("T10").Value = 0.005940112214253741

Write-Host "Done"